# Apply the "update new orleans xlsx files" edit:
#  1. Insert a new "State" column into hotel_info (between Hotel_Name and City)
#     and populate it with "Louisiana" for the existing hotel row.
#  2. Reorder the worksheet tabs so review_info comes before hotel_info.

$wb = $excel.ActiveWorkbook

$hotelWs = $wb.Worksheets.Item("hotel_info")
$reviewWs = $wb.Worksheets.Item("review_info")

# Insert a new column C (State) into hotel_info, shifting City/Zip/... right.
$hotelWs.Columns.Item(3).Insert()
$hotelWs.Cells.Item(1, 3).Value = "State"
$hotelWs.Cells.Item(2, 3).Value = "Louisiana"

# Move review_info so it becomes the first tab, ahead of hotel_info.
$reviewWs.Move($hotelWs)
